$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the newly-measured ".NET MAUI" startup time in the "Datos" table ---
$ws.Range("C4").Value = 716.9

# --- Re-colour the "Startup Time" bar-chart series (accent1 -> accent5, darker) ---
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.Format.Fill.ForeColor.ObjectThemeColor = 9   # msoThemeColorAccent5
$series.Format.Fill.ForeColor.Brightness = -0.5        # lumMod 50%

# Push the freshly entered value into the chart's cached series data as well
$series.Points().Item(1).Value = 716.9

# --- Leave the cursor where the author left it ---
$ws.Range("G20").Select()
